$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (reuse the existing header style)
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D1").Value = "238 Error"
$ws.Range("E1").Value = "206 Error"

# Build the number-format + centered-alignment style on one cell,
# then propagate it via a format-only paste so the whole range
# ends up on a single new style entry (avoids extra intermediate styles).
$template = $ws.Range("D2")
$template.NumberFormat = "0.0"
$template.HorizontalAlignment = -4108
$template.Copy()
$dataRange = $ws.Range("D2:E52")
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(2, 14.93594903635892, 13.453995975064799),
    @(3, 8.7616469498884726, 46.530318419628813),
    @(4, 35.256641446030471, 114.57384539555369),
    @(5, 19.07871864514766, 27.522078968247683),
    @(6, 8.507919439091836, 14.528292216708678),
    @(7, 64.566596060886781, 14.252632289586586),
    @(8, 46.3602693531368, 7.0152132774614415),
    @(9, 12.360170176578663, 82.715163702493783),
    @(10, 27.607062422016924, 55.366576166041085),
    @(11, 17.42649938380265, 31.026272499314189),
    @(12, 21.841914054915435, 20.125603533513697),
    @(13, 18.080632324085485, 27.9086763539147),
    @(14, 32.411394431523206, 75.924196086093957),
    @(15, 32.842823744248165, 24.933418371787582),
    @(16, 20.260840965920721, 56.709556182600267),
    @(17, 28.980553759108375, 33.429958154714825),
    @(18, 26.084745924716628, 29.497372244956807),
    @(19, 11.63220256701905, 18.408927718758662),
    @(20, 16.514900753769382, 10.302723349828057),
    @(21, 25.583380572513988, 38.820383719307472),
    @(22, 17.570731115137193, 119.98249413693733),
    @(23, 22.417394168158523, 21.455147619861634),
    @(24, 44.583515531892715, 27.556703665964562),
    @(25, 22.665061936716484, 70.068729124722353),
    @(26, 27.113497635422959, 26.45806848565519),
    @(27, 22.417885612442433, 44.35433454014867),
    @(28, 18.465086725060587, 28.902742376506978),
    @(29, 4.2097018326925308, 21.111356594467338),
    @(30, 80.066495692962349, 120.16852327267628),
    @(31, 9.6892162543084339, 7.6355784972492415),
    @(32, 17.018187146023763, 24.338915828040172),
    @(33, 10.337412312791486, 22.510396531032171),
    @(34, 35.33818526693193, 117.57962967278488),
    @(35, 6.5068846339024162, 20.702287049642962),
    @(36, 20.083814556476511, 32.419569782229985),
    @(37, 10.985369916239307, 19.750013106831091),
    @(38, 7.7845920649667164, 34.176050183839664),
    @(39, 36.668605889849459, 21.743170929549422),
    @(40, 17.332363383247582, 19.35518065479863),
    @(41, 27.602704659402661, 130.50189932743143),
    @(42, 8.8956105132135974, 31.093845484323424),
    @(43, 320.66243775594569, 81.769236345035154),
    @(44, 64.794799461847617, 175.88425218310482),
    @(45, 32.615753544294876, 35.881119645479203),
    @(46, 37.157358898641093, 56.846086489482218),
    @(47, 40.674308787380824, 20.957748179820612),
    @(48, 39.093545558640983, 13.939150526896981),
    @(49, 33.905309849780281, 32.214371413194215),
    @(50, 45.877930989128686, 15.447111629789561),
    @(51, 23.553140255949245, 46.635935648666248),
    @(52, 45.013809856838691, 9.4789595638183641),
)

foreach ($item in $data) {
    $r = $item[0]
    $d = $item[1]
    $e = $item[2]
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

$ws.Range("E2:E52").Select()
